# trackhubs.xlsx example-data refresh
# - "Hub Data": hub name now uses an underscore (Test Hub -> Test_Hub)
# - "Tracks Data": file paths no longer carry the "test-hub/" directory
#   prefix and are renumbered (test1/test2/test3), the Track2_1 related
#   specimen id is corrected, and the rich-text specimen-id cell for
#   Track1 is normalized to plain text.
# - leaves the workbook with "Hub Data" as the active sheet/selection.

$wb = $excel.ActiveWorkbook

$hubData = $wb.Worksheets.Item("Hub Data")
$hubData.Range("A2").Value = "Test_Hub"

$tracksData = $wb.Worksheets.Item("Tracks Data")
$tracksData.Range("B2").Value = "test1.Bigbed"
$tracksData.Range("F2").Value = "SAMEA104728908, SAMEA104728909 "
$tracksData.Range("B3").Value = "test2.Bigbed"
$tracksData.Range("F3").Value = "SAMEA104728909, SAMEA104728907"
$tracksData.Range("B4").Value = "test3.Bigbed"

$tracksData.Activate()
$tracksData.Range("F3").Select()

$hubData.Activate()
$hubData.Range("A2").Select()
